$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> @(newPrice, newVolume). $null price means leave D unchanged.
$updates = @{
    2  = @("51.897.35", "  +1.75%  ")
    3  = @("3.009.03", "  +3.52%  ")
    4  = @("0.999", "  +0.01%  ")
    5  = @("385.69", "  +3.86%  ")
    6  = @("105.33", "  +3.56%  ")
    7  = @("0.549", "  +1.58%  ")
    8  = @("0.999", "  -0.04%  ")
    9  = @("0.602", "  +2.74%  ")
    10 = @("37.70", "  +2.53%  ")
    11 = @($null, "  +0.30%  ")
    12 = @("0.0851", "  +2.14%  ")
    13 = @("3.480.05", "  +3.45%  ")
    14 = @("18.56", "  +1.68%  ")
    15 = @("7.67", "  +4.24%  ")
    16 = @("1.03", "  +11.70%  ")
    17 = @("2.998.68", "  +2.95%  ")
    18 = @("51.831.77", "  +1.67%  ")
    19 = @("3.33", "  +2.90%  ")
    20 = @($null, "  +4.34%  ")
    21 = @("13.07", "  +1.31%  ")
    22 = @("0.0₃0970", "  +3.13%  ")
    23 = @("69.32", "  +1.88%  ")
    24 = @("264.58", "  +2.24%  ")
    25 = @("2.94", "  +9.93%  ")
    26 = @("8.39", "  +19.43%  ")
    27 = @($null, "  +23.57%  ")
    28 = @($null, "  +2.51%  ")
    29 = @($null, "  +14.75%  ")
    30 = @("26.24", "  +2.70%  ")
    31 = @($null, "  +0.03%  ")
    32 = @("9.96", "  +1.21%  ")
    33 = @("35.31", "  +3.71%  ")
    34 = @("51.21", "  -0.10%  ")
    35 = @($null, "  +8.75%  ")
    36 = @($null, "  -1.72%  ")
    37 = @($null, "  -0.28%  ")
    38 = @("3.08", "  +3.14%  ")
    39 = @("17.27", "  +1.77%  ")
    40 = @("2.63", "  +2.14%  ")
    41 = @($null, "  +1.49%  ")
    42 = @($null, "  +4.00%  ")
    43 = @("122.24", "  +2.42%  ")
    44 = @("22.04", "  +0.94%  ")
    45 = @($null, "  +19.16%  ")
    46 = @($null, "  -1.47%  ")
    47 = @("3.34", "  +6.37%  ")
    48 = @($null, "  +2.75%  ")
    49 = @("2.049.32", "  +1.74%  ")
    50 = @("0.0337", "  +10.08%  ")
    51 = @("0.870", "  +2.80%  ")
}

# Rows whose new Price text would otherwise be auto-parsed by Excel as a
# numeric value (e.g. "0.999", "385.69"). For these we force the cell's
# number format to Text before writing the value so the literal digit
# string is preserved exactly as in the source data.
$forceTextRows = @(4,5,6,7,8,9,10,12,14,15,16,19,21,23,24,25,26,30,32,33,34,38,39,40,43,44,47,50,51)

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $newPrice = $vals[0]
    $newVolume = $vals[1]

    if ($null -ne $newPrice) {
        $dCell = $ws.Range("D$row")
        if ($forceTextRows -contains $row) {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $newPrice
    }
    $ws.Range("E$row").Value = $newVolume
}
